$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- small permutation tables pasted into columns N..U ---

# Block 1 (rows 2-3, cols N:P) - blue header row, green second row
$ws.Range("N2").Value = 2
$ws.Range("O2").Value = 5
$ws.Range("P2").Value = 6
$ws.Range("N3").Value = 6
$ws.Range("O3").Value = 5
$ws.Range("P3").Value = 2

# Block 2 (rows 5-6, cols N:Q) - blue header row, green second row
$ws.Range("N5").Value = 2
$ws.Range("O5").Value = 5
$ws.Range("P5").Value = 8
$ws.Range("Q5").Value = 6
$ws.Range("N6").Value = 6
$ws.Range("O6").Value = 8
$ws.Range("P6").Value = 5
$ws.Range("Q6").Value = 2

# Block 3 (rows 8-9, cols N:U) - blue header row, plain centered second row
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 2
$ws.Range("P8").Value = 3
$ws.Range("Q8").Value = 4
$ws.Range("R8").Value = 5
$ws.Range("S8").Value = 6
$ws.Range("T8").Value = 7
$ws.Range("U8").Value = 8

$ws.Range("N9").Value = 6
$ws.Range("O9").Value = 1
$ws.Range("P9").Value = 7
$ws.Range("Q9").Value = 5
$ws.Range("R9").Value = 4
$ws.Range("S9").Value = 2
$ws.Range("T9").Value = 8
$ws.Range("U9").Value = 3

# --- formatting: blue fill (Accent1, Lighter 40%) for header rows ---
$blue = $ws.Range("N2:P2,N5:Q5,N8:U8")
$blue.HorizontalAlignment = -4108
$blue.VerticalAlignment = -4108
$blue.WrapText = $true
$blue.Interior.ThemeColor = 5
$blue.Interior.TintAndShade = 0.6

# --- formatting: green fill (Accent3, Lighter 40%) for second rows ---
$green = $ws.Range("N3:P3,N6:Q6")
$green.HorizontalAlignment = -4108
$green.VerticalAlignment = -4108
$green.WrapText = $true
$green.Interior.ThemeColor = 7
$green.Interior.TintAndShade = 0.4

# --- plain centered style (matches pre-existing col default) for row 9 ---
$plain = $ws.Range("N9:U9")
$plain.HorizontalAlignment = -4108
$plain.VerticalAlignment = -4108

# --- view state ---
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("T12").Select()
